$d = $word.ActiveDocument

# --- TreeMap section: drop the "Use Case" bullet, the trailing blank
#     paragraph, and the whole ConcurrentHashMap sub-section (heading +
#     its two bullets). Delete from the end backwards so earlier indices
#     stay valid.
$d.Paragraphs.Item(29).Range.Delete()   # Use Case: Suitable for multi-threaded environments...
$d.Paragraphs.Item(28).Range.Delete()   # Efficiency: Similar to HashMap but designed for concurrent access...
$d.Paragraphs.Item(27).Range.Delete()   # ConcurrentHashMap:
# Paragraphs 25-26 ("Use Case: Useful..." + the trailing blank paragraph)
# are now the last two paragraphs in the story; the blank one is also the
# document's very last paragraph mark, which Range.Delete can't drop on
# its own, so remove both together in one combined range.
$p25 = $d.Paragraphs.Item(25)
$p26 = $d.Paragraphs.Item(26)
$d.Range($p25.Range.Start, $p26.Range.End).Delete()

# --- ArrayList section: widen the trailing blank paragraph's indent,
#     then remove its "Use Case" bullet.
$d.Paragraphs.Item(22).LeftIndent = 90   # 1080 -> 1800 twips (90 pt)
$d.Paragraphs.Item(21).Range.Delete()    # Use Case: Suitable for maintaining an ordered list...

# --- HashMap section: widen the trailing blank paragraph's indent, then
#     remove its "Use Case" bullet.
$d.Paragraphs.Item(18).LeftIndent = 90   # 1080 -> 1800 twips (90 pt)
$d.Paragraphs.Item(17).Range.Delete()    # Use Case: Ideal for fast lookups of products...

# --- First bullet list: drop Memory Management / Maintainability /
#     Data Integrity bullets and the first of the two trailing blank
#     paragraphs (keeping the second blank paragraph in place).
$d.Paragraphs.Item(11).Range.Delete()   # (blank ListParagraph, ind=1080)
$d.Paragraphs.Item(10).Range.Delete()   # Data Integrity: Proper data structures ensure...
$d.Paragraphs.Item(9).Range.Delete()    # Maintainability: Well-structured data makes...
$d.Paragraphs.Item(8).Range.Delete()    # Memory Management: Efficient data structures help...
